$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -18.22486370885845
$ws.Range("C2").Value = -18.22486370885845
$ws.Range("D2").Value = -18.22486370885845
$ws.Range("E2").Value = -18.22486370885845
$ws.Range("F2").Value = -18.22486370885845
$ws.Range("G2").Value = -18.22486370885845
$ws.Range("H2").Value = -18.22486370885845
$ws.Range("I2").Value = -18.22486370885845
$ws.Range("J2").Value = -18.22486370885845
$ws.Range("K2").Value = -18.22486370885845

$ws.Range("B3").Value = -18.22486370885845
$ws.Range("C3").Value = -18.22486370885845
$ws.Range("D3").Value = -18.22486370885845
$ws.Range("E3").Value = -18.22486370885845
$ws.Range("F3").Value = -18.22486370885845
$ws.Range("G3").Value = -18.22486370885845
$ws.Range("H3").Value = -18.22486370885845
$ws.Range("I3").Value = 0.4895193359593978
$ws.Range("J3").Value = -18.22486370885845
$ws.Range("K3").Value = -18.22486370885845

$ws.Range("B4").Value = -18.22486370885845
$ws.Range("C4").Value = -18.22486370885845
$ws.Range("D4").Value = -0.5213859035489609
$ws.Range("E4").Value = -18.22486370885845
$ws.Range("F4").Value = 3.934649513192435
$ws.Range("G4").Value = -18.22486370885845
$ws.Range("H4").Value = 1.719975581974993
$ws.Range("I4").Value = -18.22486370885845
$ws.Range("J4").Value = 2.583004514655447
$ws.Range("K4").Value = -18.22486370885845

$ws.Range("B5").Value = -18.22486370885845
$ws.Range("C5").Value = -18.22486370885845
$ws.Range("D5").Value = -18.22486370885845
$ws.Range("E5").Value = -18.22486370885845
$ws.Range("F5").Value = -18.22486370885845
$ws.Range("G5").Value = 3.536684683095578
$ws.Range("H5").Value = -18.22486370885845
$ws.Range("I5").Value = -18.22486370885845
$ws.Range("J5").Value = -18.22486370885845
$ws.Range("K5").Value = -18.22486370885845

$ws.Range("B6").Value = -18.22486370885845
$ws.Range("C6").Value = -18.22486370885845
$ws.Range("D6").Value = -18.22486370885845
$ws.Range("E6").Value = -18.22486370885845
$ws.Range("F6").Value = -18.22486370885845
$ws.Range("G6").Value = -18.22486370885845
$ws.Range("H6").Value = -18.22486370885845
$ws.Range("I6").Value = -18.22486370885845
$ws.Range("J6").Value = -18.22486370885845
$ws.Range("K6").Value = -18.22486370885845

$ws.Range("B7").Value = 2.745338398294459
$ws.Range("C7").Value = -18.22486370885845
$ws.Range("D7").Value = -18.22486370885845
$ws.Range("E7").Value = -18.22486370885845
$ws.Range("F7").Value = -18.22486370885845
$ws.Range("G7").Value = -18.22486370885845
$ws.Range("H7").Value = -18.22486370885845
$ws.Range("I7").Value = -18.22486370885845
$ws.Range("J7").Value = -18.22486370885845
$ws.Range("K7").Value = -18.22486370885845

$ws.Range("B8").Value = -18.22486370885845
$ws.Range("C8").Value = -18.22486370885845
$ws.Range("D8").Value = -18.22486370885845
$ws.Range("E8").Value = 2.34502990170679
$ws.Range("F8").Value = -18.22486370885845
$ws.Range("G8").Value = -18.22486370885845
$ws.Range("H8").Value = -18.22486370885845
$ws.Range("I8").Value = -18.22486370885845
$ws.Range("J8").Value = -18.22486370885845
$ws.Range("K8").Value = -18.22486370885845

$ws.Range("B9").Value = 3.732754523268161
$ws.Range("C9").Value = -18.22486370885845
$ws.Range("D9").Value = -18.22486370885845
$ws.Range("E9").Value = -18.22486370885845
$ws.Range("F9").Value = -18.22486370885845
$ws.Range("G9").Value = -18.22486370885845
$ws.Range("H9").Value = -18.22486370885845
$ws.Range("I9").Value = -18.22486370885845
$ws.Range("J9").Value = -18.22486370885845
$ws.Range("K9").Value = -18.22486370885845

$ws.Range("B10").Value = -18.22486370885845
$ws.Range("C10").Value = -18.22486370885845
$ws.Range("D10").Value = -18.22486370885845
$ws.Range("E10").Value = -18.22486370885845
$ws.Range("F10").Value = -18.22486370885845
$ws.Range("G10").Value = -18.22486370885845
$ws.Range("H10").Value = -18.22486370885845
$ws.Range("I10").Value = 0.375444873182043
$ws.Range("J10").Value = -18.22486370885845
$ws.Range("K10").Value = 2.014134551471673

$ws.Range("B11").Value = -18.22486370885845
$ws.Range("C11").Value = -18.22486370885845
$ws.Range("D11").Value = -18.22486370885845
$ws.Range("E11").Value = 1.091552627247489
$ws.Range("F11").Value = -18.22486370885845
$ws.Range("G11").Value = 1.466731589007117
$ws.Range("H11").Value = -18.22486370885845
$ws.Range("I11").Value = -18.22486370885845
$ws.Range("J11").Value = -18.22486370885845
$ws.Range("K11").Value = 1.361875310300517

$ws.Range("B12").Value = -18.22486370885845
$ws.Range("C12").Value = -18.22486370885845
$ws.Range("D12").Value = -18.22486370885845
$ws.Range("E12").Value = -18.22486370885845
$ws.Range("F12").Value = -18.22486370885845
$ws.Range("G12").Value = -18.22486370885845
$ws.Range("H12").Value = -18.22486370885845
$ws.Range("I12").Value = -18.22486370885845
$ws.Range("J12").Value = -18.22486370885845
$ws.Range("K12").Value = -18.22486370885845

$ws.Range("B13").Value = -18.22486370885845
$ws.Range("C13").Value = -18.22486370885845
$ws.Range("D13").Value = -18.22486370885845
$ws.Range("E13").Value = 1.809324026780631
$ws.Range("F13").Value = -18.22486370885845
$ws.Range("G13").Value = -18.22486370885845
$ws.Range("H13").Value = -18.22486370885845
$ws.Range("I13").Value = -18.22486370885845
$ws.Range("J13").Value = 0.9543188953739439
$ws.Range("K13").Value = 2.573162655435735

$ws.Range("B14").Value = -18.22486370885845
$ws.Range("C14").Value = -18.22486370885845
$ws.Range("D14").Value = 2.023153015564042
$ws.Range("E14").Value = -18.22486370885845
$ws.Range("F14").Value = -18.22486370885845
$ws.Range("G14").Value = -18.22486370885845
$ws.Range("H14").Value = -18.22486370885845
$ws.Range("I14").Value = -18.22486370885845
$ws.Range("J14").Value = -18.22486370885845
$ws.Range("K14").Value = 1.682087056492097

$ws.Range("B15").Value = -18.22486370885845
$ws.Range("C15").Value = -18.22486370885845
$ws.Range("D15").Value = -0.4181064274379442
$ws.Range("E15").Value = -18.22486370885845
$ws.Range("F15").Value = -18.22486370885845
$ws.Range("G15").Value = -18.22486370885845
$ws.Range("H15").Value = -18.22486370885845
$ws.Range("I15").Value = -18.22486370885845
$ws.Range("J15").Value = -18.22486370885845
$ws.Range("K15").Value = -18.22486370885845

$ws.Range("B16").Value = -18.22486370885845
$ws.Range("C16").Value = -18.22486370885845
$ws.Range("D16").Value = -18.22486370885845
$ws.Range("E16").Value = -18.22486370885845
$ws.Range("F16").Value = -18.22486370885845
$ws.Range("G16").Value = -18.22486370885845
$ws.Range("H16").Value = -18.22486370885845
$ws.Range("I16").Value = -18.22486370885845
$ws.Range("J16").Value = 2.377944782030988
$ws.Range("K16").Value = -18.22486370885845

$ws.Range("B17").Value = -18.22486370885845
$ws.Range("C17").Value = -18.22486370885845
$ws.Range("D17").Value = -1.850011008880296
$ws.Range("E17").Value = -18.22486370885845
$ws.Range("F17").Value = -18.22486370885845
$ws.Range("G17").Value = -18.22486370885845
$ws.Range("H17").Value = 2.434687159045499
$ws.Range("I17").Value = -0.02468827964941632
$ws.Range("J17").Value = 1.796665304718477
$ws.Range("K17").Value = -18.22486370885845

$ws.Range("B18").Value = -18.22486370885845
$ws.Range("C18").Value = -18.22486370885845
$ws.Range("D18").Value = -18.22486370885845
$ws.Range("E18").Value = -18.22486370885845
$ws.Range("F18").Value = -18.22486370885845
$ws.Range("G18").Value = -18.22486370885845
$ws.Range("H18").Value = 2.151669763796905
$ws.Range("I18").Value = -0.5470270951981818
$ws.Range("J18").Value = 1.764814895607314
$ws.Range("K18").Value = -18.22486370885845

$ws.Range("B19").Value = -18.22486370885845
$ws.Range("C19").Value = -18.22486370885845
$ws.Range("D19").Value = 3.149473763067693
$ws.Range("E19").Value = -18.22486370885845
$ws.Range("F19").Value = -18.22486370885845
$ws.Range("G19").Value = -18.22486370885845
$ws.Range("H19").Value = 1.646309085195373
$ws.Range("I19").Value = 1.725511350177392
$ws.Range("J19").Value = -18.22486370885845
$ws.Range("K19").Value = -18.22486370885845

$ws.Range("B20").Value = -18.22486370885845
$ws.Range("C20").Value = -18.22486370885845
$ws.Range("D20").Value = 2.416684472215509
$ws.Range("E20").Value = -18.22486370885845
$ws.Range("F20").Value = 2.235277179471344
$ws.Range("G20").Value = -18.22486370885845
$ws.Range("H20").Value = 0.9228207028889501
$ws.Range("I20").Value = 3.623434097735898
$ws.Range("J20").Value = -18.22486370885845
$ws.Range("K20").Value = 2.080771691878449

$ws.Range("B21").Value = -18.22486370885845
$ws.Range("C21").Value = 4.321923621173798
$ws.Range("D21").Value = -18.22486370885845
$ws.Range("E21").Value = 3.214657268972178
$ws.Range("F21").Value = -18.22486370885845
$ws.Range("G21").Value = 2.49336896470712
$ws.Range("H21").Value = 0.8715451903991042
$ws.Range("I21").Value = -18.22486370885845
$ws.Range("J21").Value = -18.22486370885845
$ws.Range("K21").Value = -18.22486370885845
